# Apply updated crypto price/volume figures (and a USDe/Fetch.AI row swap)
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.145.55'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '3.171.76'
$ws.Range('E3').Value = '  -4.41%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.93'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.24%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.167.70'
$ws.Range('E8').Value = '  -4.52%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  -5.47%  '
$ws.Range('E11').Value = '  -4.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.455'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.96'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = '3.688.65'
$ws.Range('E15').Value = '  -4.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.119'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '3.168.05'
$ws.Range('E17').Value = '  -4.39%  '
$ws.Range('D18').Value = '63.058.83'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.59'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -4.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '461.84'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.91'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('E22').Value = '  -4.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.65'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.47'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.43'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -3.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.76'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.76'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.04'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.60%  '
$ws.Range('E32').Value = '  -6.26%  '
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('E34').Value = '  -6.34%  '
$ws.Range('E35').Value = '  -6.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.84'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.37'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.98%  '
$ws.Range('D38').Value = '0.0₃0709'
$ws.Range('E38').Value = '  -4.33%  '
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '406.48'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.12'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.45%  '
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').Value = '2.815.77'
$ws.Range('E43').Value = '  -9.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.60'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.94%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.13'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.37%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.42'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.91'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.93'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -8.32%  '
